# Dividend Calculation workbook update
# - Adds a "2017" monthly dividend table (columns J:O) to the "Yearly" sheet,
#   mirroring the existing "2016" table (columns B:G).
# - Links the "All Time" sheet's 2017 Taxable Account figure (F8) to the new
#   Yearly!L3 cell instead of a hard-coded literal.
# - Restores the previously-active sheet/selection state.

$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# ---------------------------------------------------------------------------
# 1) Yearly sheet: build the 2017 table in columns J:O (mirrors B:G / 2016)
#    NOTE: values/formulas are written FIRST and formatting is copied in
#    afterwards - copying formats onto a cell before it holds its final
#    formula can leave that formula's cached value stale.
# ---------------------------------------------------------------------------

# Year heading (J1), mirrors B1
$wsYearly.Range("J1").Value2 = 2017

# Header row (K2:O2), mirrors C2:G2's values
$wsYearly.Range("K2").Value2 = $wsYearly.Range("C2").Value2
$wsYearly.Range("L2").Value2 = $wsYearly.Range("D2").Value2
$wsYearly.Range("M2").Value2 = $wsYearly.Range("E2").Value2
$wsYearly.Range("N2").Value2 = $wsYearly.Range("F2").Value2
$wsYearly.Range("O2").Value2 = $wsYearly.Range("G2").Value2

# Month values for the new table: Taxable Account / 401K / Suzie's Roth IRA
$months = @(
    @{Row=3;  Month=1;  Taxable=24.7; K401=7.55; Roth=0},
    @{Row=4;  Month=2;  Taxable=0;    K401=0;    Roth=0},
    @{Row=5;  Month=3;  Taxable=0;    K401=0;    Roth=0},
    @{Row=6;  Month=4;  Taxable=0;    K401=0;    Roth=0},
    @{Row=7;  Month=5;  Taxable=0;    K401=0;    Roth=0},
    @{Row=8;  Month=6;  Taxable=0;    K401=0;    Roth=0},
    @{Row=9;  Month=7;  Taxable=0;    K401=0;    Roth=0},
    @{Row=10; Month=8;  Taxable=0;    K401=0;    Roth=0},
    @{Row=11; Month=9;  Taxable=0;    K401=0;    Roth=0},
    @{Row=12; Month=10; Taxable=0;    K401=0;    Roth=0},
    @{Row=13; Month=11; Taxable=0;    K401=0;    Roth=0},
    @{Row=14; Month=12; Taxable=0;    K401=0;    Roth=0}
)

foreach ($m in $months) {
    $r = $m.Row

    # Values / formulas
    $wsYearly.Cells.Item($r, 10).Value2 = $m.Month                    # J: month number
    $wsYearly.Cells.Item($r, 11).Value2 = $wsYearly.Cells.Item($r, 3).Value2  # K: month name (mirrors C{r})
    $wsYearly.Cells.Item($r, 12).Value2 = $m.Taxable                  # L: Taxable Account
    $wsYearly.Cells.Item($r, 13).Value2 = $m.K401                     # M: 401K
    $wsYearly.Cells.Item($r, 14).Value2 = $m.Roth                     # N: Suzie's Roth IRA
    $wsYearly.Range("O" + $r).Formula = "=SUM(L" + $r + ":N" + $r + ")"  # O: Grand Total
}

# Totals row 15 (K15:O15), mirrors C15:G15
$wsYearly.Range("K15").Value2 = $wsYearly.Range("C15").Value2
$wsYearly.Range("L15").Formula = "=SUM(L3:L14)"
$wsYearly.Range("M15").Formula = "=SUM(M3:M14)"
$wsYearly.Range("N15").Formula = "=SUM(N3:N14)"
$wsYearly.Range("O15").Formula = "=SUM(O3:O14)"

# Now copy the formatting over from the mirrored 2016 columns (B:G) so that
# every new cell (J1, K2:O2, J3:O14, K15:O15) picks up the same cellXf as
# its 2016 counterpart, without disturbing the formula values set above.
$wsYearly.Range("B1").Copy()
$wsYearly.Range("J1").PasteSpecial(-4122)

$wsYearly.Range("C2:G2").Copy()
$wsYearly.Range("K2").PasteSpecial(-4122)

foreach ($m in $months) {
    $r = $m.Row
    $wsYearly.Range("B" + $r + ":G" + $r).Copy()
    $wsYearly.Range("J" + $r).PasteSpecial(-4122)
}

$wsYearly.Range("C15:G15").Copy()
$wsYearly.Range("K15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) All Time sheet: 2017 Taxable Account (F8) now pulls from Yearly!L3
# ---------------------------------------------------------------------------
$wsAllTime.Range("F8").Formula = "=Yearly!L3"

# ---------------------------------------------------------------------------
# 3) Restore view/selection state: Yearly tab active with L4 selected,
#    All Time showing L22 selected (not the active tab).
# ---------------------------------------------------------------------------
$wsAllTime.Activate()
$wsAllTime.Range("L22").Select()

$wsYearly.Activate()
$wsYearly.Range("L4").Select()
